$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextCell "D2" '29.439.92'
Set-TextCell "E2" '  +0.20%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextCell "D3" '1.848.68'
Set-TextCell "E3" '  +0.29%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextCell "D4" '1.001'
Set-TextCell "E4" '  +0.24%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextCell "D5" '240.72'
Set-TextCell "E5" '  +0.73%  '

$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextCell "D6" '0.6264'
Set-TextCell "E6" '  -0.79%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextCell "D7" '1.001'
Set-TextCell "E7" '  +0.15%  '

$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell "D8" '0.07667'
Set-TextCell "E8" '  +1.83%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell "D9" '0.2915'
Set-TextCell "E9" '  -0.40%  '

$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell "D10" '24.75'
Set-TextCell "E10" '  +1.33%  '

$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell "D11" '0.07751'
Set-TextCell "E11" '  +0.59%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell "D12" '1.851.19'
Set-TextCell "E12" '  -1.74%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell "D13" '5.028'
Set-TextCell "E13" '  +0.63%  '

$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell "D14" '0.6809'
Set-TextCell "E14" '  +0.28%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell "D15" '0.00001070'
Set-TextCell "E15" '  +3.20%  '

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell "D16" '83.49'
Set-TextCell "E16" '  +0.61%  '

$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell "D17" '6.166'
Set-TextCell "E17" '  -0.10%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell "D18" '29.468.29'
Set-TextCell "E18" '  +0.16%  '

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell "D19" '228.28'
Set-TextCell "E19" '  +0.15%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell "D20" '12.38'
Set-TextCell "E20" '  -0.22%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell "D21" '1.001'
Set-TextCell "E21" '  +0.13%  '

$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell "D22" '7.416'
Set-TextCell "E22" '  -0.49%  '

$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell "D23" '1.002'
Set-TextCell "E23" '  +0.16%  '

$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell "D24" '157.65'
Set-TextCell "E24" '  +0.48%  '

$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell "D25" '0.1372'
Set-TextCell "E25" '  -1.32%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell "D26" '8.396'
Set-TextCell "E26" '  +0.37%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell "D27" '17.71'
Set-TextCell "E27" '  +0.81%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell "D28" '1.351'
Set-TextCell "E28" '  +5.88%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell "D29" '1.465'
Set-TextCell "E29" '  +0.68%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell "D30" '0.05628'
Set-TextCell "E30" '  +0.12%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell "D31" '4.115'
Set-TextCell "E31" '  +0.28%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell "D32" '4.031'
Set-TextCell "E32" '  +0.21%  '

$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell "D33" '1.841'
Set-TextCell "E33" '  +0.38%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell "D34" '1.161'
Set-TextCell "E34" '  +0.39%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell "D35" '0.7032'
Set-TextCell "E35" '  -1.04%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell "D36" '2.596'
Set-TextCell "E36" '  +0.31%  '

$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell "D37" '1.228.99'
Set-TextCell "E37" '  -1.21%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell "D38" '2.767'
Set-TextCell "E38" '  +0.04%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell "D39" '0.01788'
Set-TextCell "E39" '  -1.08%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell "D40" '6.537'
Set-TextCell "E40" '  +3.45%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell "D41" '0.9026'
Set-TextCell "E41" '  +0.42%  '

$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell "D42" '1.002'
Set-TextCell "E42" '  +0.20%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell "D43" '101.82'
Set-TextCell "E43" '  -0.10%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell "D44" '65.86'
Set-TextCell "E44" '  +0.26%  '

$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell "D45" '0.00000000121'
Set-TextCell "E45" '  +2.42%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell "D46" '7.150'
Set-TextCell "E46" '  +1.06%  '

$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell "D47" '0.4012'
Set-TextCell "E47" '  +0.27%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell "D48" '9.006'
Set-TextCell "E48" '  +1.45%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell "D49" '0.1153'
Set-TextCell "E49" '  +2.94%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell "D50" '1.671'
Set-TextCell "E50" '  +0.04%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell "D51" '0.05715'
Set-TextCell "E51" '  +0.22%  '

